# Updates cryptos list prices/volume deltas per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.512.71"
$ws.Cells.Item(2, 5).Value = "  -2.48%  "

$ws.Cells.Item(3, 4).Value = "2.356.14"
$ws.Cells.Item(3, 5).Value = "  +0.02%  "

$ws.Cells.Item(4, 5).Value = "  -0.18%  "

$ws.Cells.Item(5, 4).Value = "'325.49"
$ws.Cells.Item(5, 5).Value = "  +3.36%  "

$ws.Cells.Item(6, 4).Value = "'100.27"
$ws.Cells.Item(6, 5).Value = "  -9.25%  "

$ws.Cells.Item(7, 5).Value = "  -1.25%  "

$ws.Cells.Item(8, 5).Value = "  -0.01%  "

$ws.Cells.Item(9, 4).Value = "'0.623"
$ws.Cells.Item(9, 5).Value = "  -2.15%  "

$ws.Cells.Item(10, 4).Value = "'39.90"
$ws.Cells.Item(10, 5).Value = "  -8.14%  "

$ws.Cells.Item(11, 5).Value = "  -1.93%  "

$ws.Cells.Item(12, 4).Value = "'8.41"
$ws.Cells.Item(12, 5).Value = "  -5.17%  "

$ws.Cells.Item(13, 4).Value = "'1.00"
$ws.Cells.Item(13, 5).Value = "  -4.03%  "

$ws.Cells.Item(14, 5).Value = "  +0.01%  "

$ws.Cells.Item(15, 5).Value = "  +0.72%  "

$ws.Cells.Item(16, 4).Value = "2.711.52"
$ws.Cells.Item(16, 5).Value = "  +0.10%  "

$ws.Cells.Item(17, 4).Value = "2.351.83"
$ws.Cells.Item(17, 5).Value = "  -3.01%  "

$ws.Cells.Item(18, 4).Value = "'8.06"
$ws.Cells.Item(18, 5).Value = "  +10.84%  "

$ws.Cells.Item(19, 4).Value = "42.624.08"
$ws.Cells.Item(19, 5).Value = "  -2.13%  "

$ws.Cells.Item(20, 5).Value = "  -2.16%  "

$ws.Cells.Item(21, 4).Value = "'76.26"
$ws.Cells.Item(21, 5).Value = "  +0.65%  "

$ws.Cells.Item(22, 5).Value = "  +7.22%  "

$ws.Cells.Item(23, 4).Value = "'265.84"
$ws.Cells.Item(23, 5).Value = "  +3.53%  "

$ws.Cells.Item(24, 5).Value = "  -10.53%  "

$ws.Cells.Item(25, 4).Value = "'10.03"
$ws.Cells.Item(25, 5).Value = "  +9.52%  "

$ws.Cells.Item(26, 4).Value = "'1.01"
$ws.Cells.Item(26, 5).Value = "  +0.51%  "

$ws.Cells.Item(27, 4).Value = "'11.45"
$ws.Cells.Item(27, 5).Value = "  -5.35%  "

$ws.Cells.Item(28, 4).Value = "'22.96"
$ws.Cells.Item(28, 5).Value = "  +2.22%  "

$ws.Cells.Item(29, 5).Value = "  -2.17%  "

$ws.Cells.Item(30, 4).Value = "'175.34"
$ws.Cells.Item(30, 5).Value = "  +0.66%  "

$ws.Cells.Item(31, 5).Value = "  -2.58%  "

$ws.Cells.Item(32, 4).Value = "'0.0898"
$ws.Cells.Item(32, 5).Value = "  -3.83%  "

$ws.Cells.Item(33, 4).Value = "'35.27"
$ws.Cells.Item(33, 5).Value = "  -10.38%  "

$ws.Cells.Item(34, 4).Value = "'6.01"
$ws.Cells.Item(34, 5).Value = "  -0.69%  "

$ws.Cells.Item(35, 5).Value = "  -0.40%  "

$ws.Cells.Item(36, 4).Value = "'4.56"
$ws.Cells.Item(36, 5).Value = "  -8.84%  "

$ws.Cells.Item(37, 5).Value = "  -5.28%  "

$ws.Cells.Item(38, 4).Value = "'2.94"
$ws.Cells.Item(38, 5).Value = "  +7.85%  "

$ws.Cells.Item(39, 5).Value = "  +1.04%  "

$ws.Cells.Item(40, 5).Value = "  -9.70%  "

$ws.Cells.Item(41, 5).Value = "  +1.11%  "

$ws.Cells.Item(42, 5).Value = "  -0.32%  "

$ws.Cells.Item(43, 4).Value = "'69.91"
$ws.Cells.Item(43, 5).Value = "  -3.99%  "

$ws.Cells.Item(44, 5).Value = "  -0.08%  "

$ws.Cells.Item(45, 4).Value = "'119.45"
$ws.Cells.Item(45, 5).Value = "  +6.91%  "

$ws.Cells.Item(46, 4).Value = "'89.98"
$ws.Cells.Item(46, 5).Value = "  +19.89%  "

$ws.Cells.Item(47, 4).Value = "'11.84"
$ws.Cells.Item(47, 5).Value = "  -8.33%  "

$ws.Cells.Item(48, 5).Value = "  -2.95%  "

$ws.Cells.Item(49, 4).Value = "'9.20"
$ws.Cells.Item(49, 5).Value = "  -1.32%  "

$ws.Cells.Item(50, 5).Value = "  -4.22%  "

$ws.Cells.Item(51, 4).Value = "1.563.40"
$ws.Cells.Item(51, 5).Value = "  +4.10%  "
